# EBEGU-1575: move the "_GoBack" bookmark from its old location (an empty
# paragraph near the end of the "Text" block) to the very start of the
# document (the first paragraph of the first cell of the address table).
#
# Word automatically renumbers every w:id when bookmarks are added/removed,
# so all we have to do is delete the old bookmark and add the new one in
# the right place; the Bookmarks collection / OOXML serializer takes care
# of reassigning ids 0..n in document order.

$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark (previously sitting in its own
#    empty paragraph further down in the letter body).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Re-create "_GoBack" collapsed at the very beginning of the document,
#    i.e. the start of the first paragraph in the first cell of the
#    address table.
$firstCell = $d.Tables.Item(1).Cell(1, 1)
$startRange = $firstCell.Range
$startRange.End = $startRange.Start

$d.Bookmarks.Add("_GoBack", $startRange)
